$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 17
$ws.Range("E17").Value = 129
$ws.Range("F17").Value = 65
$ws.Range("H17").Value = 97

# Row 18
$ws.Range("E18").Value = 121
$ws.Range("F18").Value = 54
$ws.Range("H18").Value = 90

# Row 22
$ws.Range("E22").Value = 6
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 3

# Row 29
$ws.Range("E29").Value = 18
$ws.Range("F29").Value = 11
$ws.Range("H29").Value = 14

# Row 40
$ws.Range("E40").Value = 25
$ws.Range("F40").Value = 14
$ws.Range("H40").Value = 17

# Row 41
$ws.Range("E41").Value = 44
$ws.Range("F41").Value = 19
$ws.Range("H41").Value = 30

# Row 42
$ws.Range("E42").Value = 38
$ws.Range("F42").Value = 21
$ws.Range("H42").Value = 30

# Row 50
$ws.Range("E50").Value = 25

# Row 57
$ws.Range("E57").Value = 15

# Row 65
$ws.Range("E65").Value = 35

# Row 68
$ws.Range("E68").Value = 19

# Row 72
$ws.Range("E72").Value = 46

# Row 76
$ws.Range("E76").Value = 53

# Row 87
$ws.Range("E87").Value = 18
